# SAV-1059: add a new "Invoice Product" reference-data row that is
# missing its source fields (sourceRecordId / category), to be used by
# the "invalid invoice product - missing source" importer test fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3: id, code, name, (sourceRecordId left blank), (category left blank),
# discountable, visibilityStatus
$ws.Range("A3").Value = "InvoiceProduct-cat"
$ws.Range("B3").Value = "cat"
$ws.Range("C3").Value = "Cat"
$ws.Range("F3").Value = $true
$ws.Range("F3").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("G3").Value = "current"

# Leave the selection where the editor ended up after entering the row.
$ws.Range("C10").Select()
